$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.146.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.64%  "

$ws.Range("D3").Value = "'1.910.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.29%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'251.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.5100"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.73%  "

$ws.Range("D8").Value = "'45.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.02%  "

$ws.Range("D9").Value = "'0.2955"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.69%  "

$ws.Range("D10").Value = "'0.06789"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.81%  "

$ws.Range("D11").Value = "'1.911.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.29%  "

$ws.Range("D12").Value = "'17.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("D13").Value = "'0.07356"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.50%  "

$ws.Range("D14").Value = "'0.6901"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.73%  "

$ws.Range("D15").Value = "'86.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.39%  "

$ws.Range("D16").Value = "'4.872"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.98%  "

$ws.Range("D17").Value = "'30.154.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.72%  "

$ws.Range("D18").Value = "'0.000008098"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.03%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").Value = "'12.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.32%  "

$ws.Range("D21").Value = "'2.159.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.23%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "'4.829"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.34%  "

$ws.Range("D24").Value = "'5.732"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.71%  "

$ws.Range("D25").Value = "'9.121"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.24%  "

$ws.Range("D26").Value = "'146.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.14%  "

$ws.Range("D27").Value = "'134.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.77%  "

$ws.Range("D28").Value = "'17.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.55%  "

$ws.Range("D29").Value = "'1.993"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.68%  "

$ws.Range("D30").Value = "'1.398"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").Value = "'4.230"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("D32").Value = "'0.08775"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.58%  "

$ws.Range("D33").Value = "'3.989"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.83%  "

$ws.Range("D34").Value = "'0.05069"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.89%  "

$ws.Range("D35").Value = "'1.145"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.68%  "

$ws.Range("D36").Value = "'0.7127"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.42%  "

$ws.Range("D37").Value = "'2.691"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").Value = "'2.809"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.98%  "

$ws.Range("D39").Value = "'2.280"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.63%  "

$ws.Range("D40").Value = "'0.9696"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").Value = "'0.01694"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.01%  "

$ws.Range("D42").Value = "'6.088"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.12%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'104.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.45%  "

$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").Value = "'7.600"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.74%  "

$ws.Range("D47").Value = "'0.1279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.13%  "

$ws.Range("D48").Value = "'0.05741"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("D49").Value = "'33.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.12%  "

$ws.Range("D50").Value = "'8.452"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.67%  "

$ws.Range("D51").Value = "'0.3799"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.03%  "
